# Updated symbol list on Wed Jan 18 20:24:15 UTC 2023 with GitHub Actions
# Applies refreshed Price (D) and Volume(1h) (E) values for the crypto
# ranking sheet. Values are plain text cells (inlineStr in the source
# workbook), so each write temporarily forces Text number format to stop
# Excel's automatic "looks like a number/percent" conversion, writes the
# literal string, then restores the cell style to Normal (style index 0)
# so no incidental formatting diff is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "293.02"
    "E2" = "-3.11%"
    "D3" = "30.96"
    "E3" = "-3.61%"
    "D4" = "4.890"
    "E4" = "-1.62%"
    "D5" = "0.07273"
    "E5" = "-8.17%"
    "D6" = "1.784"
    "E6" = "-16.89%"
    "D7" = "7.682"
    "E7" = "-2.21%"
    "D8" = "3.762"
    "E8" = "-1.05%"
    "D9" = "0.9025"
    "E9" = "-2.61%"
    "D10" = "0.1660"
    "D11" = "0.07573"
    "E11" = "-5.60%"
    "D12" = "0.08073"
    "E12" = "-8.28%"
    "D13" = "0.03061"
    "E13" = "-3.14%"
    "E14" = "-0.49%"
    "D15" = "0.001510"
    "E15" = "-0.03%"
    "D16" = "0.005658"
    "E16" = "-4.91%"
    "D17" = "3.460"
    "E17" = "-0.31%"
    "D18" = "2.107"
    "E18" = "-7.55%"
    "D19" = "0.3297"
    "E19" = "0.29%"
    "D20" = "0.1306"
    "E20" = "1.32%"
    "D21" = "4.366"
    "E21" = "4.57%"
    "D22" = "0.2002"
    "E22" = "11.76%"
    "D23" = "0.04484"
    "E23" = "-2.63%"
    "D24" = "0.001216"
    "E24" = "-1.68%"
    "D25" = "0.004037"
    "E25" = "-10.24%"
    "E26" = "0.16%"
    "D39" = "0.01655"
    "E39" = "-4.76%"
    "D40" = "0.04386"
    "D41" = "0.007437"
    "E41" = "0.77%"
    "D42" = "0.1319"
    "E42" = "-3.66%"
    "D43" = "0.002041"
    "E43" = "-13.52%"
    "D44" = "0.01017"
    "E44" = "-8.51%"
    "D45" = "0.00006051"
    "E45" = "0.51%"
    "E46" = "0.16%"
    "D47" = "2.172"
    "E47" = "164.66%"
    "D48" = "0.002403"
    "E48" = "-29.11%"
    "D49" = "0.00002103"
    "E49" = "0.16%"
    "D50" = "0.0002003"
    "E50" = "0.16%"
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
    $range.Style = "Normal"
}
